$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Gesellschaft" -> "Unternehmen" and "Gesellschaftskuerzel" -> "Unternehmenskuerzel"
$ws.Range("A2").Value = "Unternehmen"
$ws.Range("A3").Value = "Unternehmenskuerzel"

# Update the selection to A11
$ws.Range("A11").Select()
